# Refresh the cryptos worksheet with the latest scraped coinranking.com values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as text, even when it looks like a number
# (e.g. "1.001" or "0.6942"), by using a leading quote-prefix character,
# the same as typing '1.001 directly into Excel. Without this, Excel
# silently reinterprets such strings as numbers.
function Set-TextValue($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'$text"
}

$ws.Cells.Item(2, 4).Value = '29.122.39'
$ws.Cells.Item(2, 5).Value = '  -1.93%  '

$ws.Cells.Item(3, 4).Value = '1.852.63'
$ws.Cells.Item(3, 5).Value = '  -0.77%  '

Set-TextValue 4 4 "1.001"
$ws.Cells.Item(4, 5).Value = '  +0.19%  '

$ws.Cells.Item(5, 2).Value = 'XRP'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue 5 4 "0.6942"
$ws.Cells.Item(5, 5).Value = '  -4.94%  '

$ws.Cells.Item(6, 2).Value = 'BNB'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue 6 4 "238.97"
$ws.Cells.Item(6, 5).Value = '  -0.75%  '

$ws.Cells.Item(7, 5).Value = '  +0.16%  '

Set-TextValue 8 4 "0.07660"
$ws.Cells.Item(8, 5).Value = '  +8.07%  '

Set-TextValue 9 4 "0.3033"
$ws.Cells.Item(9, 5).Value = '  -3.08%  '

Set-TextValue 10 4 "23.40"

Set-TextValue 11 4 "0.08128"
$ws.Cells.Item(11, 5).Value = '  -1.24%  '

$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12, 4).Value = '1.852.10'
$ws.Cells.Item(12, 5).Value = '  -0.31%  '

$ws.Cells.Item(13, 2).Value = 'Polygon'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 13 4 "0.7262"
$ws.Cells.Item(13, 5).Value = '  -2.60%  '

Set-TextValue 14 4 "5.234"
$ws.Cells.Item(14, 5).Value = '  -1.44%  '

Set-TextValue 15 4 "89.09"
$ws.Cells.Item(15, 5).Value = '  -3.47%  '

$ws.Cells.Item(16, 4).Value = '29.128.41'
$ws.Cells.Item(16, 5).Value = '  -1.93%  '

Set-TextValue 17 4 "5.787"

Set-TextValue 18 4 "0.000007778"
$ws.Cells.Item(18, 5).Value = '  -0.30%  '

Set-TextValue 19 4 "13.18"
$ws.Cells.Item(19, 5).Value = '  -1.28%  '

Set-TextValue 20 4 "236.66"
$ws.Cells.Item(20, 5).Value = '  -4.67%  '

Set-TextValue 21 4 "0.9996"
$ws.Cells.Item(21, 5).Value = '  +0.23%  '

$ws.Cells.Item(22, 4).Value = '2.095.92'
$ws.Cells.Item(22, 5).Value = '  -1.00%  '

Set-TextValue 23 4 "1.001"
$ws.Cells.Item(23, 5).Value = '  +0.10%  '

Set-TextValue 24 4 "7.607"
$ws.Cells.Item(24, 5).Value = '  -1.29%  '

Set-TextValue 25 4 "9.000"
$ws.Cells.Item(25, 5).Value = '  -1.74%  '

Set-TextValue 26 4 "161.41"
$ws.Cells.Item(26, 5).Value = '  -0.83%  '

Set-TextValue 27 4 "0.1445"
$ws.Cells.Item(27, 5).Value = '  -5.72%  '

Set-TextValue 28 4 "18.07"
$ws.Cells.Item(28, 5).Value = '  -2.45%  '

Set-TextValue 29 4 "1.977"
$ws.Cells.Item(29, 5).Value = '  -1.75%  '

Set-TextValue 30 4 "1.404"
$ws.Cells.Item(30, 5).Value = '  -1.89%  '

Set-TextValue 31 4 "4.478"
$ws.Cells.Item(31, 5).Value = '  -1.35%  '

$ws.Cells.Item(32, 5).Value = '  -1.92%  '

Set-TextValue 33 4 "4.019"
$ws.Cells.Item(33, 5).Value = '  -4.10%  '

Set-TextValue 34 4 "0.05230"
$ws.Cells.Item(34, 5).Value = '  -0.80%  '

Set-TextValue 35 4 "1.189"
$ws.Cells.Item(35, 5).Value = '  -3.43%  '

Set-TextValue 36 4 "1.028"
$ws.Cells.Item(36, 5).Value = '  +3.08%  '

Set-TextValue 37 4 "0.7024"
$ws.Cells.Item(37, 5).Value = '  -6.72%  '

$ws.Cells.Item(38, 5).Value = '  -1.50%  '

Set-TextValue 39 4 "0.01856"
$ws.Cells.Item(39, 5).Value = '  -3.82%  '

$ws.Cells.Item(40, 5).Value = '  -2.27%  '

Set-TextValue 41 4 "0.9298"
$ws.Cells.Item(41, 5).Value = '  +7.46%  '

Set-TextValue 42 4 "6.043"
$ws.Cells.Item(42, 5).Value = '  +0.82%  '

$ws.Cells.Item(43, 4).Value = '1.078.68'
$ws.Cells.Item(43, 5).Value = '  +2.67%  '

Set-TextValue 44 4 "0.4269"
$ws.Cells.Item(44, 5).Value = '  -4.50%  '

Set-TextValue 45 4 "70.46"
$ws.Cells.Item(45, 5).Value = '  -0.97%  '

Set-TextValue 46 4 "1.001"
$ws.Cells.Item(46, 5).Value = '  +0.13%  '

Set-TextValue 47 4 "102.98"
$ws.Cells.Item(47, 5).Value = '  -0.82%  '

Set-TextValue 48 4 "1.782"
$ws.Cells.Item(48, 5).Value = '  -2.16%  '

$ws.Cells.Item(49, 4).Value = '1.991.83'
$ws.Cells.Item(49, 5).Value = '  -1.40%  '

Set-TextValue 50 4 "9.211"
$ws.Cells.Item(50, 5).Value = '  -3.18%  '

Set-TextValue 51 4 "7.009"
$ws.Cells.Item(51, 5).Value = '  -6.17%  '
